$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12929.876556451
$ws.Range("C2").Value = 12362.502029488
$ws.Range("E2").Value = 8761.69765487054
$ws.Range("F2").Value = 32.0141535149381

$ws.Range("B3").Value = 12426.3026737692
$ws.Range("C3").Value = 11823.8218477333
$ws.Range("E3").Value = 8450.8038151262
$ws.Range("F3").Value = 341.615235952478

$ws.Range("B4").Value = 12328.2859760581
$ws.Range("C4").Value = 11713.4451991966
$ws.Range("E4").Value = 8385.90556792144
$ws.Range("F4").Value = 334.312115296583

$ws.Range("B5").Value = 12125.3961565532
$ws.Range("C5").Value = 10787.8009746926
$ws.Range("E5").Value = 8242.91095805494
$ws.Range("F5").Value = 289.785497197814

$ws.Range("B6").Value = 4569.43278405328
$ws.Range("C6").Value = 7152.62186577625
$ws.Range("E6").Value = 7641.89788694893
$ws.Range("F6").Value = 113.277489696882

$ws.Range("B7").Value = 5146.55504097224
$ws.Range("C7").Value = 7243.61200454147
$ws.Range("E7").Value = 8472.1213330356
$ws.Range("F7").Value = 151.661389065711
